$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two bug-title cells that changed wording.
$ws.Range("B14").Value = "[GM reported] FeedService.WebAPI - CVE-2021-26701"
$ws.Range("B43").Value = "Unable to download large files from the Files grid"

# Reflect the author's on-screen scroll/selection state when the file was saved.
$ws.Application.ActiveWindow.ScrollRow = 26
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("B44").Select()
